$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tbl_aux")
Write-Host $ws.Name
Write-Host $ws.Range("A88").Value()
